# Harmonize MVP calculations and improve scoring formula
# Implements new additive formula: Score = eloGain + (winRate * 15) + (games * 0.5)

$wb = $excel.ActiveWorkbook

# --- Sheet: MVP Simulator ---
# Update H4:I8 formulas for rows 4-8
$wsSim = $wb.Worksheets.Item("MVP Simulator")
for ($r = 4; $r -le 8; $r++) {
    $wsSim.Range("H$r").Formula = "=B$r + (E$r * 15) + (D$r * 0.5)"
    $wsSim.Range("I$r").Formula = "=H$r"
}

# --- Sheet: MVP Scenarios ---
# Update static values for F/G columns (Std MVP Score / Roll MVP Score)
$wsScen = $wb.Worksheets.Item("MVP Scenarios")
$wsScen.Range("F3").Value = 56.5
$wsScen.Range("G3").Value = 56.5

$wsScen.Range("F4").Value = 43.05
$wsScen.Range("G4").Value = 43.05

$wsScen.Range("F8").Value = 65.5
$wsScen.Range("G8").Value = 65.5

$wsScen.Range("F9").Value = 75.45
$wsScen.Range("G9").Value = 75.45

$wsScen.Range("F13").Value = 93.25
$wsScen.Range("G13").Value = 93.25

$wsScen.Range("F14").Value = 57.5
$wsScen.Range("G14").Value = 57.5

# --- Sheet: Documentation ---
$wsDoc = $wb.Worksheets.Item("Documentation")

# Row 13: rename "Rolling MVP Score" -> "MVP Days Score" and update description
$wsDoc.Range("A13").Value = "MVP Days Score"
$wsDoc.Range("B13").Value = "Now harmonized with the Standard MVP Score."

# Row 21: Standard MVP Score formula -> new additive formula
$wsDoc.Range("B21").Formula = "=eloGain + (winRate * 15) + (gamesPlayed * 0.5)"

# Row 22: rename "Rolling MVP Score" -> "MVP Days Score" and update formula
$wsDoc.Range("A22").Value = "MVP Days Score"
$wsDoc.Range("B22").Formula = "=eloGain + (winRate * 15) + (gamesPlayed * 0.5)"
